# Updated cryptos list: refresh Price (D) and Volume(1h) (E) figures,
# and correct the ordering of the FirstDigitalUSD / LidoDAOToken rows (50-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as text (preserving trailing zeros / thousands dots)
# instead of silently converting them to numbers.

$ws.Range("D2").Value = "70.968.89"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").Value = "3.652.43"
$ws.Range("E3").Value = "  +5.37%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'595.88"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'195.27"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("D7").Value = "'0.647"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("D8").Value = "3.645.96"
$ws.Range("E8").Value = "  +5.37%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  +6.47%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "'57.91"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "'0.0000296"
$ws.Range("E13").Value = "  +6.02%  "
$ws.Range("D14").Value = "'9.97"
$ws.Range("E14").Value = "  +4.52%  "
$ws.Range("D15").Value = "4.240.31"
$ws.Range("E15").Value = "  +5.72%  "
$ws.Range("D16").Value = "'20.25"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").Value = "3.659.27"
$ws.Range("E17").Value = "  +5.55%  "
$ws.Range("D18").Value = "71.059.06"
$ws.Range("E18").Value = "  +5.65%  "
$ws.Range("D19").Value = "'12.77"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Value = "'1.06"
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("D22").Value = "'488.38"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").Value = "'19.11"
$ws.Range("E23").Value = "  +13.79%  "
$ws.Range("D24").Value = "'5.26"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").Value = "'4.49"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").Value = "'91.53"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("D28").Value = "'11.45"
$ws.Range("E28").Value = "  +3.91%  "
$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  +5.73%  "
$ws.Range("D30").Value = "'32.83"
$ws.Range("E30").Value = "  +4.37%  "
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  +9.01%  "
$ws.Range("D33").Value = "'12.26"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").Value = "'621.29"
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("D35").Value = "'66.26"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").Value = "'40.15"
$ws.Range("E36").Value = "  +7.00%  "
$ws.Range("D37").Value = "0.0₃0834"
$ws.Range("E37").Value = "  +9.79%  "
$ws.Range("D38").Value = "'0.412"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  +2.66%  "
$ws.Range("D42").Value = "3.334.02"
$ws.Range("E42").Value = "  +3.15%  "
$ws.Range("D43").Value = "'3.20"
$ws.Range("E43").Value = "  +9.50%  "
$ws.Range("E44").Value = "  +6.77%  "
$ws.Range("E45").Value = "  +10.65%  "
$ws.Range("D46").Value = "'0.0456"
$ws.Range("E46").Value = "  +4.95%  "
$ws.Range("D47").Value = "'9.64"
$ws.Range("E47").Value = "  +9.80%  "
$ws.Range("D48").Value = "'3.32"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("D49").Value = "'0.139"
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.23"
$ws.Range("E51").Value = "  +0.29%  "
